$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: submission date "3rd November 2023" -> "6rd November 2023"
# Only the leading digit changes; the superscript "rd" stays as its
# own run untouched, exactly like the source edit.
# ------------------------------------------------------------------
$dateAnchor = $d.Content
$found1 = $dateAnchor.Find.Execute("3rd November 2023", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if ($found1) {
    $digitRange = $d.Range($dateAnchor.Start, $dateAnchor.Start + 1)
    $digitRange.Text = "6"
}

# ------------------------------------------------------------------
# Change 2: split the run that reads
#   " to be entered from the keyboard and stored.  These results should then be displayed "
# into two runs:
#   " to be entered from the keyboard and stored."
#   "  These results should then be displayed "
# ------------------------------------------------------------------
$headAnchor = $d.Content
$found2 = $headAnchor.Find.Execute("to be entered from the keyboard and stored.", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if ($found2) {
    $splitPos = $headAnchor.End

    $tailAnchor = $d.Range($splitPos, $d.Content.End)
    $tailAnchor.Find.Execute("  These results should then be displayed ", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)

    $tailRange = $d.Range($splitPos, $tailAnchor.End)

    # Forces Word to materialise a distinct run at this boundary: toggling
    # a character property on and back off splits the run at the range
    # edges without leaving any visible formatting change behind (the
    # resulting rPr ends up identical to its sibling run again).
    $tailRange.Font.Bold = 1
    $tailRange.Font.Bold = 0
}
